$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Qs.6:92 / Qs.6:101 / Qs.6:103 / Qs.17:111" row up from row 19 to row 18,
# and clear the now-vacated "Level 5" label + values from row 19.

$ws.Range("G18").Value = "Qs.6:92"
$ws.Range("H18").Value = "Qs.6:101"
$ws.Range("I18").Value = "Qs.6:103"
$ws.Range("J18").Value = "Qs.17:111"

$ws.Range("G18:J18").Style = $ws.Range("G19:J19").Style

$ws.Range("F19").ClearContents()
$ws.Range("G19:J19").ClearContents()

# Update the active window selection / scroll position.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("K19").Select()
